$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.984.74'
$ws.Range("E2").Value = '  +2.62%  '
$ws.Range("D3").Value = '3.093.69'
$ws.Range("E3").Value = '  +5.16%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'580.63"
$ws.Range("E5").Value = '  +1.98%  '
$ws.Range("D6").Value = "'169.51"
$ws.Range("E6").Value = '  +6.42%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.089.53'
$ws.Range("E8").Value = '  +5.14%  '
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("E10").Value = '  -2.31%  '
$ws.Range("E11").Value = '  +3.93%  '
$ws.Range("E12").Value = '  +4.56%  '
$ws.Range("E13").Value = '  +2.01%  '
$ws.Range("E14").Value = '  +5.38%  '
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("D16").Value = '3.605.33'
$ws.Range("E16").Value = '  +5.15%  '
$ws.Range("D17").Value = '66.871.20'
$ws.Range("E17").Value = '  +2.37%  '
$ws.Range("D18").Value = "'7.19"
$ws.Range("E18").Value = '  +2.48%  '
$ws.Range("D19").Value = '3.092.92'
$ws.Range("E19").Value = '  +5.22%  '
$ws.Range("D20").Value = "'16.21"
$ws.Range("E20").Value = '  +4.94%  '
$ws.Range("D21").Value = "'466.47"
$ws.Range("E21").Value = '  +4.77%  '
$ws.Range("D22").Value = "'0.713"
$ws.Range("E22").Value = '  +2.63%  '
$ws.Range("D23").Value = "'7.49"
$ws.Range("E23").Value = '  +2.67%  '
$ws.Range("D24").Value = "'84.09"
$ws.Range("E24").Value = '  +2.00%  '
$ws.Range("D25").Value = "'13.13"
$ws.Range("E25").Value = '  +7.94%  '
$ws.Range("E26").Value = '  +2.71%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = "'8.03"
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("E30").Value = '  +0.57%  '
$ws.Range("D31").Value = "'2.68"
$ws.Range("E32").Value = '  +1.04%  '
$ws.Range("D33").Value = "'28.29"
$ws.Range("E33").Value = '  +3.43%  '
$ws.Range("E34").Value = '  +2.41%  '
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = "'1.01"
$ws.Range("E36").Value = '  +3.37%  '
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("D38").Value = "'47.85"
$ws.Range("E38").Value = '  +8.00%  '
$ws.Range("E39").Value = '  +5.57%  '
$ws.Range("E40").Value = '  +1.98%  '
$ws.Range("E41").Value = '  +4.18%  '
$ws.Range("E42").Value = '  +1.65%  '
$ws.Range("D43").Value = "'8.66"
$ws.Range("E43").Value = '  +1.86%  '
$ws.Range("E44").Value = '  -1.04%  '
$ws.Range("E45").Value = '  +2.09%  '
$ws.Range("D46").Value = "'381.77"
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("D47").Value = '2.789.92'
$ws.Range("E47").Value = '  +3.42%  '
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("D50").Value = "'24.79"
$ws.Range("E50").Value = '  +5.77%  '
$ws.Range("E51").Value = '  +1.19%  '
